# Update "想去人数" (want-to-go count) figures to the latest scraped values.
# Corresponds to an automated data refresh (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F8").Value = 201
$wsExhibit.Range("F14").Value = 12391
$wsExhibit.Range("F15").Value = 109
$wsExhibit.Range("F16").Value = 5484

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 116

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 116
$wsAll.Range("F10").Value = 201
$wsAll.Range("F16").Value = 12391
$wsAll.Range("F18").Value = 109
$wsAll.Range("F19").Value = 5484
